# Scen_TRA_Max_Stock.xlsx - "Modify maximum stock number"
# Target sheet: UCT1 (4th worksheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCT1")

# Update the "Year" column values for the four UC_T blocks (rows 7, 10, 13, 16)
# from 2020 to 2018.
$ws.Range("E7").Value = 2018
$ws.Range("E10").Value = 2018
$ws.Range("E13").Value = 2018
$ws.Range("E16").Value = 2018

# Update the Heavy Goods Truck maximum-stock multiplier from 1.12 to 1.15.
# H16 holds its own (non-shared) formula, while I16:AH16 form a shared
# formula group anchored at I16 - update each piece accordingly so the
# shared-formula grouping is preserved like the original file.
$ws.Range("H16").Formula = "=H29*1.15"
$ws.Range("I16:AH16").Formula = "=I29*1.15"

# Move the active selection on the sheet from H21 to F20.
$ws.Range("F20").Select() | Out-Null

$wb.Save() | Out-Null
